# Insert a brand new weekly record as the new row 275 for the
# "Hortaliza, Macroferia Regional de Talca - Brócoli" sheet.
# All existing rows from 275 downward shift down by one (to 276..347),
# preserving their original values/styles; the vacated row 275 is then
# populated with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 275:346 down to 276:347 by inserting a new blank row at 275.
$ws.Rows.Item(275).EntireRow.Insert()

# Populate the newly inserted row 275 with the new record.
$ws.Cells.Item(275, 1).Value = 5
$ws.Cells.Item(275, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(275, 3).Value = "Maule"
$ws.Cells.Item(275, 4).Value = 44736
$ws.Cells.Item(275, 5).Value = 7
$ws.Cells.Item(275, 6).Value = 100112023
$ws.Cells.Item(275, 7).Value = "Brócoli"
$ws.Cells.Item(275, 8).Value = "Sin especificar"
$ws.Cells.Item(275, 9).Value = "Primera"
$ws.Cells.Item(275, 10).Value = 5000
$ws.Cells.Item(275, 11).Value = 800
$ws.Cells.Item(275, 12).Value = 800
$ws.Cells.Item(275, 13).Value = 800
$ws.Cells.Item(275, 14).Value = "$/unidad"
$ws.Cells.Item(275, 15).Value = "Región del Maule"
$ws.Cells.Item(275, 16).Value = 800
$ws.Cells.Item(275, 17).Value = 1
$ws.Cells.Item(275, 18).Value = "Hortaliza"

# Keep the date column's number format consistent with the rest of column D.
$ws.Cells.Item(275, 4).NumberFormat = $ws.Cells.Item(276, 4).NumberFormat
